$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text: "التقييم من 5" -> "التقييم من 10"
$ws.Range("A1").Value = "التقييم من 10"

# Update evaluation scores in column A (rows 2-14)
$ws.Range("A2").Value = 8
$ws.Range("A3").Value = 8
$ws.Range("A4").Value = 5
$ws.Range("A5").Value = 7
$ws.Range("A6").Value = 7
$ws.Range("A7").Value = 7
$ws.Range("A8").Value = 8
$ws.Range("A9").Value = 5
$ws.Range("A10").Value = 6
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 5
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 7

# Update the active selection to A8
$ws.Range("A8").Select()
